$d = $word.ActiveDocument

# Locate the target paragraph: "Vestibulum eget velit..." (currently styled BodyText)
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("Vestibulum eget velit")) {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)

# Insert a brand-new (empty) paragraph immediately before it to hold the
# Mermaid source-code block; the new paragraph takes over $targetIndex and
# the "Vestibulum..." text shifts down to $targetIndex + 1.
$target.Range.InsertParagraphBefore()

$codeIndex = $targetIndex
$realIndex = $targetIndex + 1

$code = $d.Paragraphs.Item($codeIndex)
$code.Style = "Source Code"

# Each element is one "line" of the flowchart snippet; an empty string means
# a blank line (i.e. two consecutive line breaks with no text between them).
$lines = @(
    "flowchart LR",
    "",
    "A[Hard] -->|Text| B(Round)",
    "B --> C{Decision}",
    "C -->|One| D[Result 1]",
    "C -->|Two| E[Result 2]"
)

for ($li = 0; $li -lt $lines.Count; $li++) {
    if ($li -gt 0) {
        $codeP = $d.Paragraphs.Item($codeIndex)
        $insertPos = $codeP.Range.End - 1
        $breakPoint = $d.Range($insertPos, $insertPos)
        $breakPoint.InsertBreak(6)
    }

    $line = $lines[$li]
    if ($line.Length -gt 0) {
        $codeP2 = $d.Paragraphs.Item($codeIndex)
        $insertPos2 = $codeP2.Range.End - 1
        $runRange = $d.Range($insertPos2, $insertPos2)
        $runRange.InsertAfter($line)

        $codeP3 = $d.Paragraphs.Item($codeIndex)
        $styleEnd = $codeP3.Range.End - 1
        $styleStart = $styleEnd - $line.Length
        $styleRange = $d.Range($styleStart, $styleEnd)
        $styleRange.Style = "Verbatim Char"
    }
}

# The original paragraph keeps its text, but switches from BodyText to FirstParagraph
$realP = $d.Paragraphs.Item($realIndex)
$realP.Style = "First Paragraph"
